$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.169.67"
$ws.Range("E2").Value = "  +5.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.549.74"
$ws.Range("E3").Value = "  +6.29%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.51"
$ws.Range("E5").Value = "  +9.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "564.04"
$ws.Range("E6").Value = "  +6.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +4.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.540.16"
$ws.Range("E8").Value = "  +6.07%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.635"
$ws.Range("E10").Value = "  +4.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  +13.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.78"
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +6.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.36"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.100.16"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.536.98"
$ws.Range("E16").Value = "  +5.96%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.62"
$ws.Range("E17").Value = "  +6.11%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.122"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.143.90"
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  +8.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +4.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.05"
$ws.Range("E22").Value = "  +13.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  +11.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.41"
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.18"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.09"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.94"
$ws.Range("E27").Value = "  +9.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.31"
$ws.Range("E28").Value = "  +8.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.05"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.09"
$ws.Range("E30").Value = "  +10.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.49"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "633.73"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.66"
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.77"
$ws.Range("E34").Value = "  +5.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.112"
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.46"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0826"
$ws.Range("E37").Value = "  +14.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  +19.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.46"
$ws.Range("E39").Value = "  +5.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.389"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  +10.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.142.23"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  +10.52%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("E47").Value = "  +11.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0419"
$ws.Range("E48").Value = "  +5.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.133"
$ws.Range("E50").Value = "  +6.16%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.63"
$ws.Range("E51").Value = "  +9.01%  "
